{"js": "// Add three new character styles (GaNStyle, GaNParagraph, GaNLinks) and\n// apply them to the runs that were newly styled in the authored edit.\n\nconst doc = context.document;\n\n// --- Create the three character styles -------------------------------\ndoc.addStyle(\"GaNStyle\", Word.StyleType.character);\ndoc.addStyle(\"GaNParagraph\", Word.StyleType.character);\ndoc.addStyle(\"GaNLinks\", Word.StyleType.character);\nawait context.sync();\n\nconst gaNStyle = doc.getStyles().getByName(\"GaNStyle\");\ngaNStyle.font.name = \"Calibri\";\ngaNStyle.font.size = 14; // w:sz 28 half-points\n\nconst gaNParagraph = doc.getStyles().getByName(\"GaNParagraph\");\ngaNParagraph.font.name = \"Calibri\";\ngaNParagraph.font.size = 10; // w:sz 20 half-points\n\nconst gaNLinks = doc.getStyles().getByName(\"GaNLinks\");\ngaNLinks.font.name = \"Calibri\";\ngaNLinks.font.size = 9.5; // w:sz 19 half-points\ngaNLinks.font.bold = true;\ngaNLinks.font.color = \"#000080\";\ngaNLinks.font.underline = Word.UnderlineType.single;\nawait context.sync();\n\n// --- Apply GaNStyle to every \"Kampagnendaten 2022 ...\" run -----------\nconst kampagnendatenText =\n  \"Kampagnendaten 2022 f\u00fcr das Sternbild Herkules: 13.-22. Juni, 12.-21. Juli, 10.-19. August\";\nconst kampagnendatenResults = doc.body.search(kampagnendatenText, { matchCase: true });\nkampagnendatenResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < kampagnendatenResults.items.length; i++) {\n  kampagnendatenResults.items[i].style = \"GaNStyle\";\n}\n\n// --- Apply GaNParagraph to the \"Mach mit an einer weltweiten...\" run -\nconst machMitText =\n  \"Mach mit an einer weltweiten Kampagne, die schw\u00e4chsten sichtbaren Sterne zu beobachten und aufzuzeichnen, um die Lichtverschmutzung an einem Ort zu messen. Durch das Auffinden und Beobachten des Sternbild Herkules am Nachthimmel und den Vergleich mit den Helligkeitskarten, lernen Menschen auf der ganzen Erde, wie die Lichter in ihrer Gemeinde zur Lichtverschmutzung beitragen. Dein Beitrag zur Online-Datenbank beschreibt den sichtbaren Nachthimmel.\";\nconst machMitResults = doc.body.search(machMitText, { matchCase: true });\nmachMitResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < machMitResults.items.length; i++) {\n  machMitResults.items[i].style = \"GaNParagraph\";\n}\n\n// --- Apply GaNLinks to the \"Die Schaubilder in diesem Dokument...\" run\nconst schaubilderText =\n  \"Die Schaubilder in diesem Dokument wurden von Jan Hollan, CzechGlobe, bereitgestellt. (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\";\nconst schaubilderResults = doc.body.search(schaubilderText, { matchCase: true });\nschaubilderResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < schaubilderResults.items.length; i++) {\n  schaubilderResults.items[i].style = \"GaNLinks\";\n}\n\nawait context.sync();\n", "ps1": "# Add three new character styles (GaNStyle, GaNParagraph, GaNLinks) and\n# apply them to the runs that were newly styled in the authored edit.\n\n$d = $word.ActiveDocument\n\n# --- Create the three character styles --------------------------------\n# wdStyleTypeCharacter = 2\n$gaNStyle = $d.Styles.Add(\"GaNStyle\", 2)\n$gaNStyle.Font.Name = \"Calibri\"\n$gaNStyle.Font.Size = 14   # w:sz 28 half-points\n\n$gaNParagraph = $d.Styles.Add(\"GaNParagraph\", 2)\n$gaNParagraph.Font.Name = \"Calibri\"\n$gaNParagraph.Font.Size = 10   # w:sz 20 half-points\n\n$gaNLinks = $d.Styles.Add(\"GaNLinks\", 2)\n$gaNLinks.Font.Name = \"Calibri\"\n$gaNLinks.Font.Size = 9.5   # w:sz 19 half-points\n$gaNLinks.Font.Bold = $true\n$gaNLinks.Font.Color = 0x800000  # wdColorNavy -> OOXML w:color 000080 (BGR-encoded \"long\" color)\n$gaNLinks.Font.Underline = 1  # wdUnderlineSingle\n\n# --- Apply GaNStyle to every \"Kampagnendaten 2022 ...\" run ------------\n$kampagnendatenText = \"Kampagnendaten 2022 f\u00fcr das Sternbild Herkules: 13.-22. Juni, 12.-21. Juli, 10.-19. August\"\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = $kampagnendatenText\n$find.MatchCase = $true\n$find.MatchWildcards = $false\nwhile ($find.Execute()) {\n    $range.Style = \"GaNStyle\"\n    $range.Collapse(0)\n}\n\n# --- Apply GaNParagraph to the \"Mach mit an einer weltweiten...\" run --\n$machMitText = \"Mach mit an einer weltweiten Kampagne, die schw\u00e4chsten sichtbaren Sterne zu beobachten und aufzuzeichnen, um die Lichtverschmutzung an einem Ort zu messen. Durch das Auffinden und Beobachten des Sternbild Herkules am Nachthimmel und den Vergleich mit den Helligkeitskarten, lernen Menschen auf der ganzen Erde, wie die Lichter in ihrer Gemeinde zur Lichtverschmutzung beitragen. Dein Beitrag zur Online-Datenbank beschreibt den sichtbaren Nachthimmel.\"\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.ClearFormatting()\n$find2.Text = $machMitText\n$find2.MatchCase = $true\n$find2.MatchWildcards = $false\nwhile ($find2.Execute()) {\n    $range2.Style = \"GaNParagraph\"\n    $range2.Collapse(0)\n}\n\n# --- Apply GaNLinks to the \"Die Schaubilder in diesem Dokument...\" run\n$schaubilderText = \"Die Schaubilder in diesem Dokument wurden von Jan Hollan, CzechGlobe, bereitgestellt. (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\"\n$range3 = $d.Content\n$find3 = $range3.Find\n$find3.ClearFormatting()\n$find3.Text = $schaubilderText\n$find3.MatchCase = $true\n$find3.MatchWildcards = $false\nwhile ($find3.Execute()) {\n    $range3.Style = \"GaNLinks\"\n    $range3.Collapse(0)\n}\n\nWrite-Output \"done\"\n"}
